$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D2").Value = "65.427.30"
$ws.Range("D3").Value = "3.553.31"
$ws.Range("D5").Value = "600.79"
$ws.Range("D6").Value = "139.95"
$ws.Range("D7").Value = "3.552.59"
$ws.Range("D11").Value = "7.01"
$ws.Range("D13").Value = "4.162.12"
$ws.Range("D15").Value = "27.11"
$ws.Range("D16").Value = "3.561.60"
$ws.Range("D18").Value = "65.491.04"
$ws.Range("D19").Value = "10.24"
$ws.Range("D20").Value = "5.89"
$ws.Range("D21").Value = "14.28"
$ws.Range("D22").Value = "395.64"
$ws.Range("D23").Value = "0.572"
$ws.Range("D24").Value = "3.700.17"
$ws.Range("D25").Value = "74.22"
$ws.Range("D26").Value = "0.999"
$ws.Range("D31").Value = "8.32"
$ws.Range("D32").Value = "3.570.67"
$ws.Range("D34").Value = "23.88"
$ws.Range("D35").Value = "0.147"
$ws.Range("D37").Value = "7.08"
$ws.Range("D39").Value = "167.35"
$ws.Range("D40").Value = "5.05"
$ws.Range("D43").Value = "26.75"
$ws.Range("D44").Value = "42.92"
$ws.Range("D47").Value = "1.69"
$ws.Range("D48").Value = "1.19"
$ws.Range("D49").Value = "2.448.72"
$ws.Range("D51").Value = "2.37"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("E11").Value = "  -6.39%  "
$ws.Range("E12").Value = "  +3.79%  "
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("E19").Value = "  +3.15%  "
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +3.99%  "
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +8.75%  "
$ws.Range("E28").Value = "  +9.87%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +4.10%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  +4.27%  "
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("E40").Value = "  +4.65%  "
$ws.Range("E41").Value = "  +3.40%  "
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  +15.23%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E47").Value = "  +3.45%  "
$ws.Range("E48").Value = "  +8.16%  "
$ws.Range("E49").Value = "  +10.19%  "
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("E51").Value = "  +20.05%  "

# --- Rows 34/35 swap (Kaspa <-> EthereumClassic) additional B/C updates ---
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"

Write-Host "Edit script completed"
